# Update countries & provincias Spain
# Applies the COVID-19 dataset refresh described by the commit:
#  - bumps the "Datos actualizados" timestamp
#  - refreshes a handful of per-country case counters
#  - re-sorts three countries (Estado de Palestina/Bielorrusia and
#    Macao/Guam/Kenia) to keep the sheet alphabetically consistent,
#    carrying each country's own row of numbers along with its name

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 27 de Marzo de 2020 a las 03:12"

# --- Plain numeric refreshes (no reordering) ---------------------------
# Suiza
$ws.Range("E11").Value = 11488
$ws.Range("G11").Value = 39
$ws.Range("H11").Value = 192

# Noruega
$ws.Range("B20").Value = 3372
$ws.Range("C20").Value = 288
$ws.Range("E20").Value = 3352

# Hong Kong
$ws.Range("B61").Value = 454
$ws.Range("C61").Value = 43
$ws.Range("E61").Value = 340

# Venezuela
$ws.Range("E98").Value = 91
$ws.Range("G98").Value = 1
$ws.Range("H98").Value = 1

# --- Estado de Palestina / Bielorrusia swap (rows 104-105) -------------
# Row 104 becomes "Estado de Palestina" with refreshed counts.
$ws.Range("A104").Value = "Estado de Palestina"
$ws.Range("B104").Value = 86
$ws.Range("C104").Value = 15
$ws.Range("D104").Value = 17
$ws.Range("E104").Value = 68
$ws.Range("F104").Value = 0
$ws.Range("G104").Value = 0
$ws.Range("H104").Value = 1

# Row 105 becomes "Bielorrusia" (its former row-104 figures, unchanged).
$ws.Range("A105").Value = "Bielorrusia"
$ws.Range("B105").Value = 86
$ws.Range("C105").Value = 0
$ws.Range("D105").Value = 29
$ws.Range("E105").Value = 57
$ws.Range("F105").Value = 2
$ws.Range("G105").Value = 0
$ws.Range("H105").Value = 0

# --- Macao / Guam / Kenia reorder (rows 128-130) ------------------------
# Row 128 becomes "Macao" with refreshed counts.
$ws.Range("A128").Value = "Macao"
$ws.Range("B128").Value = 33
$ws.Range("C128").Value = 2
$ws.Range("D128").Value = 10
$ws.Range("E128").Value = 23
$ws.Range("F128").Value = 0
$ws.Range("G128").Value = 0
$ws.Range("H128").Value = 0

# Row 129 becomes "Guam" (its former row-128 figures, unchanged).
$ws.Range("A129").Value = "Guam"
$ws.Range("B129").Value = 32
$ws.Range("C129").Value = 0
$ws.Range("D129").Value = 0
$ws.Range("E129").Value = 31
$ws.Range("F129").Value = 0
$ws.Range("G129").Value = 0
$ws.Range("H129").Value = 1

# Row 130 becomes "Kenia" (its former row-129 figures, unchanged).
$ws.Range("A130").Value = "Kenia"
$ws.Range("B130").Value = 31
$ws.Range("C130").Value = 3
$ws.Range("D130").Value = 1
$ws.Range("E130").Value = 29
$ws.Range("F130").Value = 0
$ws.Range("G130").Value = 1
$ws.Range("H130").Value = 1
